$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the other header cells (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows 2-29: I = 1 (constant), J = copy of H
for ($r = 2; $r -le 29; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
